# Scheduled market-price refresh for the Twintania_Profits leve-crafting workbook.
# Updates currentAveragePrice/NQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ columns
# (H:N) per-row on each job sheet to the latest pulled prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 129.92308
$ws.Range("I9").Value = 138
$ws.Range("J9").Value = 103
$ws.Range("K9").Value = 138
$ws.Range("L9").Value = 103
$ws.Range("M9").Value = 31
$ws.Range("N9").Value = -441
# row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 335.33334
$ws.Range("I28").Value = 314.25
$ws.Range("K28").Value = 314.25
$ws.Range("M28").Value = 170.75
# row 29 (Leve Item ID 4575)
$ws.Range("H29").Value = 1575
$ws.Range("I29").Value = 1575
$ws.Range("K29").Value = 4725
$ws.Range("M29").Value = -4444
# row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 785.7273
$ws.Range("J41").Value = 1141.4
$ws.Range("L41").Value = 1141.4
$ws.Range("N41").Value = -2021.4
# row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 41199.6
$ws.Range("I70").Value = 41199.6
$ws.Range("K70").Value = 123598.8
$ws.Range("M70").Value = -123328.8
# row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 41199.6
$ws.Range("I73").Value = 41199.6
$ws.Range("K73").Value = 123598.8
$ws.Range("M73").Value = -122662.8
# row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 1852.55
$ws.Range("I86").Value = 1873.7646
$ws.Range("K86").Value = 1873.7646
$ws.Range("M86").Value = -750.7646
# row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 1852.55
$ws.Range("I89").Value = 1873.7646
$ws.Range("K89").Value = 9368.823
$ws.Range("M89").Value = -3752.823
# row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 537.7857
$ws.Range("I92").Value = 549.9167
$ws.Range("K92").Value = 549.9167
$ws.Range("M92").Value = 698.0833
# row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1539.1052
$ws.Range("I98").Value = 1484.375
$ws.Range("K98").Value = 1484.375
$ws.Range("M98").Value = 13.625
# row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 27823.75
$ws.Range("J106").Value = 34899.75
$ws.Range("L106").Value = 34899.75
$ws.Range("N106").Value = -36161.75
# row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 1260.5483
$ws.Range("I107").Value = 979.1111
$ws.Range("K107").Value = 979.1111
$ws.Range("M107").Value = 940.8889
# row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1539.1052
$ws.Range("I122").Value = 1484.375
$ws.Range("K122").Value = 4453.125
$ws.Range("M122").Value = -2003.125
# row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 10290.85
$ws.Range("I137").Value = 7284.3184
$ws.Range("K137").Value = 21852.9552
$ws.Range("M137").Value = -19302.9552

$ws = $wb.Worksheets.Item("ARM")
# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 3617.5813
$ws.Range("I2").Value = 2680.5588
$ws.Range("K2").Value = 2680.5588
$ws.Range("M2").Value = -2567.5588
# row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 2509.84
$ws.Range("I32").Value = 1549.3636
$ws.Range("K32").Value = 1549.3636
$ws.Range("M32").Value = -1262.3636
# row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 6984.5
$ws.Range("I45").Value = 9356
$ws.Range("K45").Value = 9356
$ws.Range("M45").Value = -8979
# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 3617.5813
$ws.Range("I116").Value = 2680.5588
$ws.Range("K116").Value = 2680.5588
$ws.Range("M116").Value = -386.5587999999998

$ws = $wb.Worksheets.Item("BSM")
# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 3617.5813
$ws.Range("I3").Value = 2680.5588
$ws.Range("K3").Value = 2680.5588
$ws.Range("M3").Value = -2566.5588
# row 88 (Leve Item ID 10626)
$ws.Range("H88").Value = 21614.166
$ws.Range("J88").Value = 21614.166
$ws.Range("L88").Value = 21614.166
$ws.Range("N88").Value = -22426.166
# row 91 (Leve Item ID 10626)
$ws.Range("H91").Value = 21614.166
$ws.Range("J91").Value = 21614.166
$ws.Range("L91").Value = 21614.166
$ws.Range("N91").Value = -24422.166

$ws = $wb.Worksheets.Item("CRP")
# row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 254.25
$ws.Range("I22").Value = 256.17648
$ws.Range("J22").Value = 243.33333
$ws.Range("K22").Value = 256.17648
$ws.Range("L22").Value = 243.33333
$ws.Range("M22").Value = 93.82351999999997
$ws.Range("N22").Value = -943.3333299999999
# row 103 (Leve Item ID 19558)
$ws.Range("H103").Value = 24361.375
$ws.Range("I103").Value = 24361.375
$ws.Range("K103").Value = 24361.375
$ws.Range("M103").Value = -23189.375

$ws = $wb.Worksheets.Item("CUL")
# row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 533.9091
$ws.Range("I2").Value = 229.8
$ws.Range("K2").Value = 1378.8
$ws.Range("M2").Value = -1265.8
# row 75 (Leve Item ID 12863)
$ws.Range("H75").Value = 483.42856
$ws.Range("J75").Value = 219.5
$ws.Range("L75").Value = 658.5
$ws.Range("N75").Value = -2654.5
# row 78 (Leve Item ID 12863)
$ws.Range("H78").Value = 483.42856
$ws.Range("J78").Value = 219.5
$ws.Range("L78").Value = 1975.5
$ws.Range("N78").Value = -11959.5
# row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 1162.3334
$ws.Range("I122").Value = 750
$ws.Range("J122").Value = 1368.5
$ws.Range("K122").Value = 6750
$ws.Range("L122").Value = 12316.5
$ws.Range("M122").Value = -4300
$ws.Range("N122").Value = -17216.5
# row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 2038.0714
$ws.Range("I132").Value = 1497.5
$ws.Range("K132").Value = 13477.5
$ws.Range("M132").Value = -10947.5

$ws = $wb.Worksheets.Item("GSM")
# row 5 (Leve Item ID 1681)
$ws.Range("H5").Value = 11740.6
$ws.Range("I5").Value = 11175.75
$ws.Range("K5").Value = 11175.75
$ws.Range("M5").Value = -11063.75
# row 18 (Leve Item ID 4309)
$ws.Range("H18").Value = 50000
$ws.Range("I18").Value = 50000
$ws.Range("K18").Value = 50000
$ws.Range("M18").Value = -49707
# row 40 (Leve Item ID 4113)
$ws.Range("H40").Value = 28499.75
$ws.Range("I40").Value = 27999.666
$ws.Range("K40").Value = 27999.666
$ws.Range("M40").Value = -27848.666
# row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 23196.727
$ws.Range("J43").Value = 35566.43
$ws.Range("L43").Value = 35566.43
$ws.Range("N43").Value = -35868.43
# row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 38889.668
$ws.Range("I46").Value = 13347
$ws.Range("J46").Value = 51661
$ws.Range("K46").Value = 13347
$ws.Range("L46").Value = 51661
$ws.Range("M46").Value = -13191
$ws.Range("N46").Value = -51973
# row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 34950
$ws.Range("J57").Value = 36261.875
$ws.Range("L57").Value = 36261.875
$ws.Range("N57").Value = -37901.875

$ws = $wb.Worksheets.Item("LTW")
# row 24 (Leve Item ID 3774)
$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20686

$ws = $wb.Worksheets.Item("WVR")
# row 10 (Leve Item ID 3313)
$ws.Range("H10").Value = 49999.332
$ws.Range("I10").Value = 49999.332
$ws.Range("K10").Value = 49999.332
$ws.Range("M10").Value = -49830.332
# row 29 (Leve Item ID 3568)
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4248
# row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -21240
# row 80 (Leve Item ID 10911)
$ws.Range("H80").Value = 34833
$ws.Range("J80").Value = 37249.5
$ws.Range("L80").Value = 37249.5
$ws.Range("N80").Value = -39245.5
# row 83 (Leve Item ID 10911)
$ws.Range("H83").Value = 34833
$ws.Range("J83").Value = 37249.5
$ws.Range("L83").Value = 111748.5
$ws.Range("N83").Value = -121732.5
# row 88 (Leve Item ID 10842)
$ws.Range("H88").Value = 23250
$ws.Range("J88").Value = 23250
$ws.Range("L88").Value = 23250
$ws.Range("N88").Value = -24062
# row 91 (Leve Item ID 10842)
$ws.Range("H91").Value = 23250
$ws.Range("J91").Value = 23250
$ws.Range("L91").Value = 23250
$ws.Range("N91").Value = -26058
# row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 8211.857
$ws.Range("I126").Value = 9196.6
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 27589.8
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -25119.8
$ws.Range("N126").Value = -22190
